$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read current "Tipo" header text (currently in D1) and "single" value (currently in D2)
$tipoHeader = $ws.Range("D1").Text
$tipoValue = $ws.Range("D2").Text

# Insert new "MAE" header at D1 (pushing "Tipo" header to E1)
$ws.Range("D1").Value = "MAE"
$ws.Range("E1").Value = $tipoHeader

# Copy the header style from D1 (bold/border/center) onto the new E1 cell
$ws.Range("D1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Move "single" value to E2, set new MAE numeric value in D2
$ws.Range("E2").Value = $tipoValue
$ws.Range("D2").Value = 0.607111138630155

# Update existing B2 / C2 values
$ws.Range("B2").Value = 0.5482355203125806
$ws.Range("C2").Value = 0.9967343174295523
